$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 6129.3335
$ws.Range("I43").Value = 2001
$ws.Range("J43").Value = 8193.5
$ws.Range("K43").Value = 2001
$ws.Range("L43").Value = 8193.5
$ws.Range("M43").Value = -1932
$ws.Range("N43").Value = -8331.5

$ws.Range("H74").Value = 5333.1665
$ws.Range("I74").Value = 5333.1665
$ws.Range("K74").Value = 5333.1665
$ws.Range("M74").Value = -4397.1665

$ws.Range("H77").Value = 5333.1665
$ws.Range("I77").Value = 5333.1665
$ws.Range("K77").Value = 26665.8325
$ws.Range("M77").Value = -21985.8325

$ws.Range("H99").Value = 281.8
$ws.Range("I99").Value = 136.33333
$ws.Range("K99").Value = 408.99999
$ws.Range("M99").Value = 1089.00001

$ws.Range("H106").Value = 28251.666
$ws.Range("I106").Value = 28251.666
$ws.Range("K106").Value = 28251.666
$ws.Range("M106").Value = -27620.666

$ws.Range("H135").Value = 1994
$ws.Range("I135").Value = 1399.6666
$ws.Range("K135").Value = 12596.9994
$ws.Range("M135").Value = -10061.9994

$ws.Range("H137").Value = 3455.0715
$ws.Range("I137").Value = 2124.25
$ws.Range("J137").Value = 5229.5
$ws.Range("K137").Value = 6372.75
$ws.Range("L137").Value = 15688.5
$ws.Range("M137").Value = -3822.75
$ws.Range("N137").Value = -20788.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7702.1
$ws.Range("I32").Value = 5659.7427
$ws.Range("K32").Value = 5659.7427
$ws.Range("M32").Value = -5372.7427

$ws.Range("H132").Value = 3037.9
$ws.Range("I132").Value = 2819.889
$ws.Range("K132").Value = 8459.667000000001
$ws.Range("M132").Value = -5929.667000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 687.5
$ws.Range("I22").Value = 687.5
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 687.5
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -514.5
$ws.Range("N22").ClearContents()

$ws.Range("H82").Value = 99999
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 99999
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 99999
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -100765

$ws.Range("H85").Value = 99999
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 99999
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 99999
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -102651

$ws.Range("H86").Value = 1828.5555
$ws.Range("J86").Value = 1442.4
$ws.Range("L86").Value = 1442.4
$ws.Range("N86").Value = -3688.4

$ws.Range("H89").Value = 1828.5555
$ws.Range("J89").Value = 1442.4
$ws.Range("L89").Value = 7212
$ws.Range("N89").Value = -18444

$ws.Range("H134").Value = 1666.1
$ws.Range("I134").Value = 1437.2632
$ws.Range("K134").Value = 4311.7896
$ws.Range("M134").Value = -1776.7896

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5087.154
$ws.Range("I31").Value = 3431.625
$ws.Range("J31").Value = 7736
$ws.Range("K31").Value = 3431.625
$ws.Range("L31").Value = 7736
$ws.Range("M31").Value = -3136.625
$ws.Range("N31").Value = -8326

$ws.Range("H34").Value = 5087.154
$ws.Range("I34").Value = 3431.625
$ws.Range("J34").Value = 7736
$ws.Range("K34").Value = 3431.625
$ws.Range("L34").Value = 7736
$ws.Range("M34").Value = -3229.625
$ws.Range("N34").Value = -8140

$ws.Range("H58").Value = 3000.1428
$ws.Range("I58").Value = 991.6667
$ws.Range("J58").Value = 4506.5
$ws.Range("K58").Value = 991.6667
$ws.Range("L58").Value = 4506.5
$ws.Range("M58").Value = -788.6667
$ws.Range("N58").Value = -4912.5

$ws.Range("H132").Value = 2941.375
$ws.Range("I132").Value = 1584
$ws.Range("J132").Value = 7013.5
$ws.Range("K132").Value = 4752
$ws.Range("L132").Value = 21040.5
$ws.Range("M132").Value = -2222
$ws.Range("N132").Value = -26100.5

$ws.Range("H136").Value = 3000.1428
$ws.Range("I136").Value = 991.6667
$ws.Range("J136").Value = 4506.5
$ws.Range("K136").Value = 2975.0001
$ws.Range("L136").Value = 13519.5
$ws.Range("M136").Value = -425.0001000000002
$ws.Range("N136").Value = -18619.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2002
$ws.Range("I80").Value = 2002
$ws.Range("K80").Value = 6006
$ws.Range("M80").Value = -5070

$ws.Range("H83").Value = 2002
$ws.Range("I83").Value = 2002
$ws.Range("K83").Value = 18018
$ws.Range("M83").Value = -13338

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 553.8333
$ws.Range("I2").Value = 64.8
$ws.Range("K2").Value = 64.8
$ws.Range("M2").Value = 48.2

$ws.Range("H126").Value = 4781.1
$ws.Range("I126").Value = 4506
$ws.Range("K126").Value = 13518
$ws.Range("M126").Value = -11048

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3284.0833
$ws.Range("I16").Value = 1790.9
$ws.Range("J16").Value = 10750
$ws.Range("K16").Value = 1790.9
$ws.Range("L16").Value = 10750
$ws.Range("M16").Value = -1620.9
$ws.Range("N16").Value = -11090

$ws.Range("H82").Value = 5357.875
$ws.Range("J82").Value = 5199.6
$ws.Range("L82").Value = 5199.6
$ws.Range("N82").Value = -5921.6

$ws.Range("H85").Value = 5357.875
$ws.Range("J85").Value = 5199.6
$ws.Range("L85").Value = 5199.6
$ws.Range("N85").Value = -7695.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 10642.777
$ws.Range("J81").Value = 12196.5
$ws.Range("L81").Value = 24393
$ws.Range("N81").Value = -26515

$ws.Range("H84").Value = 10642.777
$ws.Range("J84").Value = 12196.5
$ws.Range("L84").Value = 121965
$ws.Range("N84").Value = -132573

$ws.Range("H100").Value = 1972.75
$ws.Range("I100").Value = 1997
$ws.Range("K100").Value = 3994
$ws.Range("M100").Value = -3453

$ws.Range("H132").Value = 1067.6666
$ws.Range("I132").Value = 1067.6666
$ws.Range("K132").Value = 3202.9998
$ws.Range("M132").Value = -672.9998000000001

$ws.Range("H136").Value = 2818.0833
$ws.Range("I136").Value = 1104.1111
$ws.Range("K136").Value = 3312.3333
$ws.Range("M136").Value = -762.3333000000002
